$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 475
$ws.Range("I135").Value = 475
$ws.Range("K135").Value = 4275
$ws.Range("M135").Value = -1740
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 1995
$ws.Range("J29").Value = 1995
$ws.Range("L29").Value = 1995
$ws.Range("N29").Value = -2611
$ws.Range("H32").Value = 2804
$ws.Range("I32").Value = 2804
$ws.Range("K32").Value = 2804
$ws.Range("M32").Value = -2517
$ws.Range("H110").Value = 5277.6665
$ws.Range("I110").Value = 5499.8
$ws.Range("K110").Value = 5499.8
$ws.Range("M110").Value = -3454.8
$ws.Range("H124").Value = 37571
$ws.Range("J124").Value = 37571
$ws.Range("L124").Value = 37571
$ws.Range("N124").Value = -47391
$ws.Range("H125").Value = 39999
$ws.Range("J125").Value = 39999
$ws.Range("L125").Value = 39999
$ws.Range("N125").Value = -49839
$ws.Range("H132").Value = 1546.0769
$ws.Range("I132").Value = 1410
$ws.Range("J132").Value = 1852.25
$ws.Range("K132").Value = 4230
$ws.Range("L132").Value = 5556.75
$ws.Range("M132").Value = -1700
$ws.Range("N132").Value = -10616.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 20000
$ws.Range("K82").Value = 20000
$ws.Range("M82").Value = -19617
$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 20000
$ws.Range("K85").Value = 20000
$ws.Range("M85").Value = -18674
$ws.Range("H94").Value = 1347.5
$ws.Range("I94").Value = 875.5
$ws.Range("J94").Value = 2999.5
$ws.Range("K94").Value = 875.5
$ws.Range("L94").Value = 2999.5
$ws.Range("M94").Value = -424.5
$ws.Range("N94").Value = -3901.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 299.5
$ws.Range("I16").Value = 299.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 299.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -12.5
$ws.Range("H31").Value = 3436
$ws.Range("J31").Value = 5123.75
$ws.Range("L31").Value = 5123.75
$ws.Range("N31").Value = -5713.75
$ws.Range("H34").Value = 3436
$ws.Range("J34").Value = 5123.75
$ws.Range("L34").Value = 5123.75
$ws.Range("N34").Value = -5527.75
$ws.Range("H69").Value = 3293.3333
$ws.Range("I69").Value = 3293.3333
$ws.Range("K69").Value = 3293.3333
$ws.Range("M69").Value = -2544.3333
$ws.Range("H72").Value = 3293.3333
$ws.Range("I72").Value = 3293.3333
$ws.Range("K72").Value = 9879.999899999999
$ws.Range("M72").Value = -6135.999899999999
$ws.Range("H88").Value = 23874.75
$ws.Range("J88").Value = 25833
$ws.Range("L88").Value = 25833
$ws.Range("N88").Value = -26645
$ws.Range("H91").Value = 23874.75
$ws.Range("J91").Value = 25833
$ws.Range("L91").Value = 25833
$ws.Range("N91").Value = -28641
$ws.Range("H107").Value = 807.3333
$ws.Range("I107").Value = 572.6667
$ws.Range("K107").Value = 572.6667
$ws.Range("M107").Value = 1347.3333
$ws.Range("H113").Value = 299.5
$ws.Range("I113").Value = 299.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 299.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1870.5
$ws.Range("H132").Value = 4602.1113
$ws.Range("J132").Value = 4977.25
$ws.Range("L132").Value = 14931.75
$ws.Range("N132").Value = -19991.75
$ws.Range("N16").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 7997.5
$ws.Range("I80").Value = 7997.5
$ws.Range("K80").Value = 23992.5
$ws.Range("M80").Value = -23056.5
$ws.Range("H83").Value = 7997.5
$ws.Range("I83").Value = 7997.5
$ws.Range("K83").Value = 71977.5
$ws.Range("M83").Value = -67297.5
$ws.Range("H107").Value = 1549
$ws.Range("J107").Value = 1747
$ws.Range("L107").Value = 5241
$ws.Range("N107").Value = -9081
$ws.Range("H109").Value = 72.333336
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H118").Value = 500
$ws.Range("I118").Value = 500
$ws.Range("K118").Value = 1500
$ws.Range("M118").Value = -257
$ws.Range("H121").Value = 1006.8571
$ws.Range("I121").Value = 253.75
$ws.Range("K121").Value = 761.25
$ws.Range("M121").Value = 548.75
$ws.Range("H122").Value = 800
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100
$ws.Range("H125").Value = 266.5
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 33
$ws.Range("K125").Value = 1500
$ws.Range("L125").Value = 99
$ws.Range("M125").Value = 3420
$ws.Range("N125").Value = -9939
$ws.Range("H137").Value = 3825.6667
$ws.Range("J137").Value = 3738.75
$ws.Range("L137").Value = 11216.25
$ws.Range("N137").Value = -21416.25
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2825.4546
$ws.Range("I97").Value = 2731.2778
$ws.Range("J97").Value = 3249.25
$ws.Range("K97").Value = 2731.2778
$ws.Range("L97").Value = 3249.25
$ws.Range("M97").Value = -2235.2778
$ws.Range("N97").Value = -4241.25
$ws.Range("H132").Value = 1686.8462
$ws.Range("I132").Value = 1660.75
$ws.Range("K132").Value = 4982.25
$ws.Range("M132").Value = -2452.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H132").Value = 1398.8
$ws.Range("I132").Value = 998
$ws.Range("K132").Value = 2994
$ws.Range("M132").Value = -464
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 934.3
$ws.Range("I136").Value = 919.125
$ws.Range("J136").Value = 995
$ws.Range("K136").Value = 2757.375
$ws.Range("L136").Value = 2985
$ws.Range("M136").Value = -207.375
$ws.Range("N136").Value = -8085
